# Leave Card update - "Upload Leave Card 12/27/2023 4:01 PM"
# Adds the period dates for rows 80-117 (PERIOD column, monthly 1st-of-month
# entries running Jan-2023 .. Feb-2026) and records 1.25 "EARNED" credits for
# the eleven already-elapsed periods (rows 80-90). Also restores the sheet
# view/selection state left behind by the editing session.

$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Sheet1")

# --- PERIOD (column A) date stamps for rows 80..117 -------------------
$periodDates = @{
    80  = 44927   # 2023-01-01
    81  = 44958   # 2023-02-01
    82  = 44986   # 2023-03-01
    83  = 45017   # 2023-04-01
    84  = 45047   # 2023-05-01
    85  = 45078   # 2023-06-01
    86  = 45108   # 2023-07-01
    87  = 45139   # 2023-08-01
    88  = 45170   # 2023-09-01
    89  = 45200   # 2023-10-01
    90  = 45231   # 2023-11-01
    91  = 45261   # 2023-12-01
    92  = 45292   # 2024-01-01
    93  = 45323   # 2024-02-01
    94  = 45352   # 2024-03-01
    95  = 45383   # 2024-04-01
    96  = 45413   # 2024-05-01
    97  = 45444   # 2024-06-01
    98  = 45474   # 2024-07-01
    99  = 45505   # 2024-08-01
    100 = 45536   # 2024-09-01
    101 = 45566   # 2024-10-01
    102 = 45597   # 2024-11-01
    103 = 45627   # 2024-12-01
    104 = 45658   # 2025-01-01
    105 = 45689   # 2025-02-01
    106 = 45717   # 2025-03-01
    107 = 45748   # 2025-04-01
    108 = 45778   # 2025-05-01
    109 = 45809   # 2025-06-01
    110 = 45839   # 2025-07-01
    111 = 45870   # 2025-08-01
    112 = 45901   # 2025-09-01
    113 = 45931   # 2025-10-01
    114 = 45962   # 2025-11-01
    115 = 45992   # 2025-12-01
    116 = 46023   # 2026-01-01
    117 = 46054   # 2026-02-01
}

foreach ($row in $periodDates.Keys) {
    $ws1.Cells.Item($row, 1).Value2 = $periodDates[$row]
}

# --- EARNED (column C) credit of 1.25 for the eleven completed periods -
for ($row = 80; $row -le 90; $row++) {
    $ws1.Cells.Item($row, 3).Value2 = 1.25
}

# --- Restore the view state captured when the file was last saved -----
# Top pane keeps its existing selection (B4:C4); scroll/select the split
# (bottomLeft) pane down to the newly-edited rows, then leave CONVERTION
# as the active tab, matching the saved workbook view.
$ws1.Activate()
$win = $excel.ActiveWindow
$win.SplitRow = 69
$ws1.Range("E93").Select()

$ws2 = $wb.Worksheets.Item("CONVERTION")
$ws2.Activate()
